$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 125, pushing the current
# rows 125-138 down to 127-140 (their contents/styles move with them).
$ws.Rows.Item(125).Insert()
$ws.Rows.Item(125).Insert()

# New row 125: Brócoli "Primera" entry for Región Metropolitana, 2021-08-19 (serial 44449)
$ws.Cells.Item(125,1).Value = 11
$ws.Cells.Item(125,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(125,3).Value = "Bíobío"
$ws.Cells.Item(125,4).Value = 44449
$ws.Cells.Item(125,5).Value = 8
$ws.Cells.Item(125,6).Value = 100112023
$ws.Cells.Item(125,7).Value = "Brócoli"
$ws.Cells.Item(125,8).Value = "Sin especificar"
$ws.Cells.Item(125,9).Value = "Primera"
$ws.Cells.Item(125,10).Value = 1000
$ws.Cells.Item(125,11).Value = 700
$ws.Cells.Item(125,12).Value = 800
$ws.Cells.Item(125,13).Value = 750
$ws.Cells.Item(125,14).Value = "`$/unidad"
$ws.Cells.Item(125,15).Value = "Región Metropolitana"
$ws.Cells.Item(125,16).Value = 750
$ws.Cells.Item(125,17).Value = 1
$ws.Cells.Item(125,18).Value = "Hortaliza"

# New row 126: Brócoli "Segunda" entry for Región Metropolitana, 2021-08-19 (serial 44449)
$ws.Cells.Item(126,1).Value = 11
$ws.Cells.Item(126,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(126,3).Value = "Bíobío"
$ws.Cells.Item(126,4).Value = 44449
$ws.Cells.Item(126,5).Value = 8
$ws.Cells.Item(126,6).Value = 100112023
$ws.Cells.Item(126,7).Value = "Brócoli"
$ws.Cells.Item(126,8).Value = "Sin especificar"
$ws.Cells.Item(126,9).Value = "Segunda"
$ws.Cells.Item(126,10).Value = 500
$ws.Cells.Item(126,11).Value = 600
$ws.Cells.Item(126,12).Value = 600
$ws.Cells.Item(126,13).Value = 600
$ws.Cells.Item(126,14).Value = "`$/unidad"
$ws.Cells.Item(126,15).Value = "Región Metropolitana"
$ws.Cells.Item(126,16).Value = 600
$ws.Cells.Item(126,17).Value = 1
$ws.Cells.Item(126,18).Value = "Hortaliza"

Write-Output "Inserted rows 125 and 126; dimension now $($ws.UsedRange.Address())"
